# PLANILLA DE EVALUACION FINAL FASE 2 - apply "Equipo 7" roster edits
# Target sheet: EVALUACION1 (the active/evaluation worksheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVALUACION1")

# --- Update team member names (B4, B5) ---
# These feed C24 (=$B$4) and C36 (=B5) plus the relevance score lookups.
$ws.Range("B4").Value = "MORA PALMA MATIAS FRANCISCO"
$ws.Range("B5").Value = "LEYTON CISTERNA SEBASTIAN ANDRES"

# --- Update rubric self-assessment selections (column C, rows 13-19) ---
# Row 14 ("Genera evidencias ...") : Logrado -> Logro incipiente
$ws.Range("C14").Value = "Logro incipiente"
# Row 16 ("Utiliza de manera precisa el lenguaje tecnico ...") : Completamente logrado -> Logrado
$ws.Range("C16").Value = "Logrado"
# Row 18 ("Entrega la documentacion y evidencias ...") : Completamente logrado -> Logrado
$ws.Range("C18").Value = "Logrado"
# Row 19 ("Generan evidencias claras dentro del repositorio ...") : Logrado -> Logro incipiente
$ws.Range("C19").Value = "Logro incipiente"

# --- Update the active selection on the sheet to match the saved view ---
$ws.Activate() | Out-Null
$ws.Range("E22").Select() | Out-Null
